$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.9
$ws.Range("I2").Value = 2.45
$ws.Range("J2").Value = 3.5
$ws.Range("L2").Value = 3.1
$ws.Range("AC2").Value = 26
$ws.Range("AD2").Value = 41
$ws.Range("AG2").Value = 17
$ws.Range("AI2").Value = 700
$ws.Range("AK2").Value = 11
$ws.Range("AL2").Value = 10
$ws.Range("AN2").Value = 21

# Row 3 updates
$ws.Range("G3").Value = 2.63
$ws.Range("I3").Value = 3
$ws.Range("L3").Value = 3.75
$ws.Range("AC3").Value = 26
$ws.Range("AK3").Value = 13
$ws.Range("AL3").Value = 13
